$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("C5").Value = 100000
$ws.Range("E5").Value = 299

# Row 6
$ws.Range("E6").Value = 281

# Row 7
$ws.Range("E7").Value = 123

# Row 8
$ws.Range("E8").Value = 9

# Row 9
$ws.Range("E9").Value = 119

# Row 10
$ws.Range("E10").Value = 85

# Row 11
$ws.Range("E11").Value = 27

# Row 12
$ws.Range("E12").Value = 12

# Row 13
$ws.Range("E13").Value = 46

# Row 14
$ws.Range("E14").Value = 3

# Row 15
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 289

# Row 17
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 289

# Row 18
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = -286

# Row 19
$ws.Range("E19").Value = 3.2

# Row 20
$ws.Range("E20").Value = 30.25

# Row 21
$ws.Range("E21").Value = 1.07

# Row 22
$ws.Range("E22").Value = 9633.33

# Row 23
$ws.Range("E23").Value = -9533.33

# Row 25
$ws.Range("C25").Value = 0
$ws.Range("E25").Value = 3.25

# Row 26
$ws.Range("C26").Value = 0
$ws.Range("E26").Value = 3.25

# Row 27
$ws.Range("C27").Value = 0
$ws.Range("E27").Value = 3.25
